# Append new daily OHLCV rows (662-671) to the Bitcoin price history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @{ Row = 662; A = 45131.41666666666; B = 30091.5;  C = 30102.82; D = 28874;    E = 29189.1;  F = 21487.26359591 },
    @{ Row = 663; A = 45132.41666666666; B = 29188.72; C = 29378.38; D = 29065.48; E = 29236.54; F = 9918.33334093 },
    @{ Row = 664; A = 45133.41666666666; B = 29235.78; C = 29690.93; D = 29111.01; E = 29358.33; F = 13931.08744133 },
    @{ Row = 665; A = 45134.41666666666; B = 29358;    C = 29572.18; D = 29087.68; E = 29225.73; F = 10636.17094649 },
    @{ Row = 666; A = 45135.41666666666; B = 29224.23; C = 29535.39; D = 29129.53; E = 29328.09; F = 9764.89835609 },
    @{ Row = 667; A = 45136.41666666666; B = 29327.86; C = 29412.87; D = 29267.84; E = 29369.53; F = 4075.55845071 },
    @{ Row = 668; A = 45137.41666666666; B = 29369.79; C = 29456.1;  D = 29053.2;  E = 29288.27; F = 6120.52992371 },
    @{ Row = 669; A = 45138.41666666666; B = 29288.05; C = 29524.5;  D = 29122.39; E = 29240.57; F = 11591.19804062 },
    @{ Row = 670; A = 45139.41666666666; B = 29237.8;  C = 29722.98; D = 28612;    E = 29705.37; F = 24547.8278278 },
    @{ Row = 671; A = 45140.41666666666; B = 29706.04; C = 30033.5;  D = 28930.67; E = 29169.38; F = 24016.90813323 }
)

# The last existing data row (661) is the formatting template: column A has a
# bold/centered/bordered custom datetime style (style index 2 in the source
# file), columns B:F are plain number cells with no special formatting.
$templateRange = $ws.Range("A661:F661")

foreach ($r in $newRows) {
    $row = $r.Row

    # Clone formatting from the template row (covers both the styled date
    # cell in column A and the plain cells in B:F) before writing values.
    $templateRange.Copy() | Out-Null
    $ws.Range("A" + $row + ":F" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}

$excel.CutCopyMode = 0
